$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("network_weights")

$ws.Range("P2").Value = -0.023514515608917513
$ws.Range("T3").Value = -0.0042985429780012276
$ws.Range("U3").Value = -0.3862455233618894
$ws.Range("H4").Value = 0.39937697108766507
$ws.Range("Q4").Value = -0.03542421905660554
$ws.Range("L5").Value = 0.15412923517892946
$ws.Range("J6").Value = -0.5642745921423303
$ws.Range("B7").Value = 0.5757567365376377
$ws.Range("Q7").Value = -2.4842397030421965
$ws.Range("J8").Value = 0.8208743930445419
$ws.Range("S8").Value = 1.1327684770336084
$ws.Range("U8").Value = 0.05955395591848739
$ws.Range("I9").Value = 0.8397556475862318
$ws.Range("V9").Value = -0.14286868649613396
$ws.Range("J10").Value = -0.6909613222544959
$ws.Range("J11").Value = -0.22265277936169836
$ws.Range("K11").Value = -0.8996511802531113
$ws.Range("F12").Value = 0.006191739051416804
$ws.Range("B13").Value = -0.7324436557326487
$ws.Range("O13").Value = -0.15618188320482398
$ws.Range("P13").Value = -0.42315373918803745
$ws.Range("V13").Value = 0.5420485402062647
$ws.Range("J14").Value = 1.1928081309101652
$ws.Range("R15").Value = -1.032691309658765
$ws.Range("V15").Value = -0.1323464725248332
$ws.Range("B16").Value = -0.4163656598160977
$ws.Range("C16").Value = -0.07616408334273689
$ws.Range("I16").Value = 0.528567328841388
$ws.Range("M16").Value = -0.635005905530747
$ws.Range("O16").Value = 0.07564813471975738
$ws.Range("P16").Value = -0.14781774867848457
